# fix(gui) step 1 and 2
# Update the price list date and the D-column prices (step 1 / step 2 pricing)
# on the "Hoja1" worksheet of the PARKER price list workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step: bump the price-list date (A1) forward by one day
$ws.Range("A1").Value = 45309

# Step: update unit prices in column D to the new pricing
$ws.Range("D24").Value = 1170.84
$ws.Range("D25").Value = 1209.705
$ws.Range("D26").Value = 1243.706
$ws.Range("D27").Value = 1327.753
$ws.Range("D28").Value = 1616.336
$ws.Range("D29").Value = 1894.714
$ws.Range("D30").Value = 2166.772
$ws.Range("D31").Value = 2341.67
$ws.Range("D36").Value = 1498.77
$ws.Range("D37").Value = 1645.97
$ws.Range("D38").Value = 1743.129
$ws.Range("D39").Value = 1848.562
$ws.Range("D40").Value = 2064.753
$ws.Range("D41").Value = 2681.749
$ws.Range("D42").Value = 2710.896
$ws.Range("D43").Value = 2778.914
$ws.Range("D44").Value = 3546.519
$ws.Range("D45").Value = 4163.515
$ws.Range("D46").Value = 5101.156
$ws.Range("D47").Value = 5839.596
$ws.Range("D52").Value = 1598.362
$ws.Range("D53").Value = 1923.861
$ws.Range("D54").Value = 2113.336
$ws.Range("D55").Value = 2302.81
$ws.Range("D56").Value = 2501.996
$ws.Range("D57").Value = 2778.914
$ws.Range("D58").Value = 2987.821
$ws.Range("D59").Value = 3240.451
$ws.Range("D60").Value = 4789.252
$ws.Range("D61").Value = 5183.748
$ws.Range("D62").Value = 5766.733
$ws.Range("D63").Value = 6918.14
$ws.Range("D68").Value = 1919.981
$ws.Range("D69").Value = 2207.581
$ws.Range("D70").Value = 2302.81
$ws.Range("D71").Value = 2501.996
$ws.Range("D72").Value = 2848.392
$ws.Range("D73").Value = 3080.122
$ws.Range("D74").Value = 3376.478
$ws.Range("D75").Value = 3751.527
$ws.Range("D76").Value = 5072.003
$ws.Range("D77").Value = 5946.487
$ws.Range("D78").Value = 6961.853
$ws.Range("D79").Value = 7505.974
$ws.Range("D85").Value = 2684.175
$ws.Range("D86").Value = 3118.994
$ws.Range("D87").Value = 3479.474
$ws.Range("D88").Value = 3833.151
$ws.Range("D89").Value = 4595.898
$ws.Range("D90").Value = 4810.629
$ws.Range("D91").Value = 5674.423
$ws.Range("D92").Value = 6364.291
$ws.Range("D93").Value = 7530.27
$ws.Range("D94").Value = 8497.058000000001
$ws.Range("D95").Value = 10382.055
$ws.Range("D101").Value = 4061.489
$ws.Range("D102").Value = 4595.898
$ws.Range("D103").Value = 5188.597
$ws.Range("D104").Value = 6150.535
$ws.Range("D105").Value = 6218.547
$ws.Range("D106").Value = 7855.767
$ws.Range("D107").Value = 8841.995999999999
$ws.Range("D108").Value = 9983.691000000001
$ws.Range("D109").Value = 11562.608
$ws.Range("D115").Value = 5827.463
$ws.Range("D116").Value = 6442.029
$ws.Range("D117").Value = 7481.691
$ws.Range("D118").Value = 7768.324
$ws.Range("D119").Value = 8744.825999999999
$ws.Range("D120").Value = 10192.585
$ws.Range("D121").Value = 11902.687
$ws.Range("D122").Value = 13457.321
$ws.Range("D123").Value = 14623.299
